# Add invalid login test cases below the existing valid username/password rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: correct username, wrong (numeric) password
$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = 1234

# Row 4: wrong username, correct password
$ws.Range("A4").Value = "saranya"
$ws.Range("B4").Value = "admin123"

# Row 5: username only, password left blank
$ws.Range("A5").Value = "Admin"

# Row 7 (row 6 intentionally skipped): password only, username left blank
$ws.Range("B7").Value = "admin123"

# Leave the selection where Excel would land after these edits
$ws.Range("C9").Select() | Out-Null
